# Auto-generated Excel COM-interop script
# Applies Universalis market-price refresh data to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1042.7142
$ws.Range("I33").Value = 1042.7142
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1042.7142
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -813.7141999999999

$ws.Range("H98").Value = 1754.9286
$ws.Range("I98").Value = 1736.1538
$ws.Range("J98").Value = 1999
$ws.Range("K98").Value = 1736.1538
$ws.Range("L98").Value = 1999
$ws.Range("M98").Value = -238.1538
$ws.Range("N98").Value = -4995

$ws.Range("H103").Value = 359.53845
$ws.Range("I103").Value = 370.1111
$ws.Range("J103").Value = 335.75
$ws.Range("K103").Value = 1110.3333
$ws.Range("L103").Value = 1007.25
$ws.Range("M103").Value = -524.3333
$ws.Range("N103").Value = -2179.25

$ws.Range("H113").Value = 100000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 100000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 100000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -106508

$ws.Range("H122").Value = 1754.9286
$ws.Range("I122").Value = 1736.1538
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5208.4614
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -2758.4614
$ws.Range("N122").Value = -10897

$ws.Range("H125").Value = 2130.3333
$ws.Range("I125").Value = 5450.5
$ws.Range("J125").Value = 470.25
$ws.Range("K125").Value = 49054.5
$ws.Range("L125").Value = 4232.25
$ws.Range("M125").Value = -46594.5
$ws.Range("N125").Value = -9152.25

$ws.Range("H132").Value = 5932.9355
$ws.Range("I132").Value = 4917.16
$ws.Range("J132").Value = 10165.333
$ws.Range("K132").Value = 14751.48
$ws.Range("L132").Value = 30495.999
$ws.Range("M132").Value = -12221.48
$ws.Range("N132").Value = -35555.999

$ws.Range("H137").Value = 12090.353
$ws.Range("I137").Value = 12618.4
$ws.Range("J137").Value = 11870.333
$ws.Range("K137").Value = 37855.2
$ws.Range("L137").Value = 35610.999
$ws.Range("M137").Value = -35305.2
$ws.Range("N137").Value = -40710.999

$ws.Range("H138").Value = 4545.5
$ws.Range("I138").Value = 1348.125
$ws.Range("J138").Value = 5529.3076
$ws.Range("K138").Value = 4044.375
$ws.Range("L138").Value = 16587.9228
$ws.Range("M138").Value = 1095.625
$ws.Range("N138").Value = -26867.9228

$ws.Range("H141").Value = 3581
$ws.Range("I141").Value = 3581
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10743
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5563

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 64250
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 64250
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 64250
$ws.Range("N43").Value = -64876

$ws.Range("H45").Value = 4777.154
$ws.Range("I45").Value = 3979.1428
$ws.Range("J45").Value = 5708.1665
$ws.Range("K45").Value = 3979.1428
$ws.Range("L45").Value = 5708.1665
$ws.Range("M45").Value = -3602.1428
$ws.Range("N45").Value = -6462.1665

$ws.Range("H122").Value = 7216.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7216.5
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 21649.5
$ws.Range("N122").Value = -26549.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3831.6667
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 4697.5
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 4697.5
$ws.Range("M22").Value = -1927
$ws.Range("N22").Value = -5043.5

$ws.Range("H99").Value = 3584.2222
$ws.Range("I99").Value = 3109.6
$ws.Range("J99").Value = 4177.5
$ws.Range("K99").Value = 3109.6
$ws.Range("L99").Value = 4177.5
$ws.Range("M99").Value = -1611.6
$ws.Range("N99").Value = -7173.5

$ws.Range("H134").Value = 15197.333
$ws.Range("I134").Value = 11287.429
$ws.Range("J134").Value = 28882
$ws.Range("K134").Value = 33862.287
$ws.Range("L134").Value = 86646
$ws.Range("M134").Value = -31327.287
$ws.Range("N134").Value = -91716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7310.4287
$ws.Range("I16").Value = 8137.5
$ws.Range("J16").Value = 6979.6
$ws.Range("K16").Value = 8137.5
$ws.Range("L16").Value = 6979.6
$ws.Range("M16").Value = -7850.5
$ws.Range("N16").Value = -7553.6

$ws.Range("H31").Value = 2534.4707
$ws.Range("I31").Value = 1292
$ws.Range("J31").Value = 3555.0715
$ws.Range("K31").Value = 1292
$ws.Range("L31").Value = 3555.0715
$ws.Range("M31").Value = -997
$ws.Range("N31").Value = -4145.0715

$ws.Range("H34").Value = 2534.4707
$ws.Range("I34").Value = 1292
$ws.Range("J34").Value = 3555.0715
$ws.Range("K34").Value = 1292
$ws.Range("L34").Value = 3555.0715
$ws.Range("M34").Value = -1090
$ws.Range("N34").Value = -3959.0715

$ws.Range("H63").Value = 60271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 60271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 60271
$ws.Range("N63").Value = -61643

$ws.Range("H66").Value = 60271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 60271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 180813
$ws.Range("N66").Value = -187677

$ws.Range("H99").Value = 10664.333
$ws.Range("I99").Value = 10994
$ws.Range("J99").Value = 10499.5
$ws.Range("K99").Value = 10994
$ws.Range("L99").Value = 10499.5
$ws.Range("M99").Value = -9496
$ws.Range("N99").Value = -13495.5

$ws.Range("H105").Value = 4692.8696
$ws.Range("I105").Value = 4069.4
$ws.Range("J105").Value = 5861.875
$ws.Range("K105").Value = 4069.4
$ws.Range("L105").Value = 5861.875
$ws.Range("M105").Value = -2322.4
$ws.Range("N105").Value = -9355.875

$ws.Range("H107").Value = 1929.6666
$ws.Range("I107").Value = 1200
$ws.Range("J107").Value = 2075.6
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 2075.6
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -5915.6

$ws.Range("H113").Value = 7310.4287
$ws.Range("I113").Value = 8137.5
$ws.Range("J113").Value = 6979.6
$ws.Range("K113").Value = 8137.5
$ws.Range("L113").Value = 6979.6
$ws.Range("M113").Value = -5967.5
$ws.Range("N113").Value = -11319.6

$ws.Range("H122").Value = 4662.154
$ws.Range("I122").Value = 4556.4
$ws.Range("J122").Value = 4728.25
$ws.Range("K122").Value = 13669.2
$ws.Range("L122").Value = 14184.75
$ws.Range("M122").Value = -11219.2
$ws.Range("N122").Value = -19084.75

$ws.Range("H126").Value = 10664.333
$ws.Range("I126").Value = 10994
$ws.Range("J126").Value = 10499.5
$ws.Range("K126").Value = 32982
$ws.Range("L126").Value = 31498.5
$ws.Range("M126").Value = -30512
$ws.Range("N126").Value = -36438.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 434
$ws.Range("I2").Value = 972.25
$ws.Range("J2").Value = 102.76923
$ws.Range("K2").Value = 5833.5
$ws.Range("L2").Value = 616.61538
$ws.Range("M2").Value = -5720.5
$ws.Range("N2").Value = -842.61538

$ws.Range("H4").Value = 8556808
$ws.Range("I4").Value = 9166994
$ws.Range("J4").Value = 1234567
$ws.Range("K4").Value = 27500982
$ws.Range("L4").Value = 3703701
$ws.Range("M4").Value = -27500870
$ws.Range("N4").Value = -3703925

$ws.Range("H14").Value = 234.55556
$ws.Range("I14").Value = 234.55556
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 703.66668
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -530.66668

$ws.Range("H93").Value = 1180
$ws.Range("I93").Value = 590.9091
$ws.Range("J93").Value = 2800
$ws.Range("K93").Value = 1772.7273
$ws.Range("L93").Value = 8400
$ws.Range("M93").Value = 99.27269999999999
$ws.Range("N93").Value = -12144

$ws.Range("H113").Value = 2647611.8
$ws.Range("I113").Value = 6666933.5
$ws.Range("J113").Value = 455254.47
$ws.Range("K113").Value = 20000800.5
$ws.Range("L113").Value = 1365763.41
$ws.Range("M113").Value = -19998630.5
$ws.Range("N113").Value = -1370103.41

$ws.Range("H139").Value = 1942.3636
$ws.Range("I139").Value = 1358.25
$ws.Range("J139").Value = 3500
$ws.Range("K139").Value = 4074.75
$ws.Range("L139").Value = 10500
$ws.Range("M139").Value = 1065.25
$ws.Range("N139").Value = -20780

$ws.Range("H141").Value = 5246.5
$ws.Range("I141").Value = 5246.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15739.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -10559.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 53333.332
$ws.Range("I104").Value = 90000
$ws.Range("J104").Value = 35000
$ws.Range("K104").Value = 90000
$ws.Range("L104").Value = 35000
$ws.Range("M104").Value = -86506
$ws.Range("N104").Value = -41988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5229.6387
$ws.Range("I7").Value = 5817.7617
$ws.Range("J7").Value = 4406.2666
$ws.Range("K7").Value = 5817.7617
$ws.Range("L7").Value = 4406.2666
$ws.Range("M7").Value = -5705.7617
$ws.Range("N7").Value = -4630.2666

$ws.Range("H22").Value = 815.2308
$ws.Range("I22").Value = 466.55554
$ws.Range("J22").Value = 1599.75
$ws.Range("K22").Value = 466.55554
$ws.Range("L22").Value = 1599.75
$ws.Range("M22").Value = -171.55554
$ws.Range("N22").Value = -2189.75

$ws.Range("H27").Value = 815.2308
$ws.Range("I27").Value = 466.55554
$ws.Range("J27").Value = 1599.75
$ws.Range("K27").Value = 466.55554
$ws.Range("L27").Value = 1599.75
$ws.Range("M27").Value = -359.55554
$ws.Range("N27").Value = -1813.75

$ws.Range("H126").Value = 5229.6387
$ws.Range("I126").Value = 5817.7617
$ws.Range("J126").Value = 4406.2666
$ws.Range("K126").Value = 17453.2851
$ws.Range("L126").Value = 13218.7998
$ws.Range("M126").Value = -14983.2851
$ws.Range("N126").Value = -18158.7998

$ws.Range("H132").Value = 2296.3948
$ws.Range("I132").Value = 1950.9
$ws.Range("J132").Value = 3592
$ws.Range("K132").Value = 5852.700000000001
$ws.Range("L132").Value = 10776
$ws.Range("M132").Value = -3322.700000000001
$ws.Range("N132").Value = -15836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2938.1667
$ws.Range("I122").Value = 3157.5
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 9472.5
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -7022.5
$ws.Range("N122").Value = -12398.5

$ws.Range("H132").Value = 12510.625
$ws.Range("I132").Value = 8251.235000000001
$ws.Range("J132").Value = 22854.857
$ws.Range("K132").Value = 24753.705
$ws.Range("L132").Value = 68564.571
$ws.Range("M132").Value = -22223.705
$ws.Range("N132").Value = -73624.571
